$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 2 (pushes existing data rows 2..12 down to 3..13).
# This works inside the ListObject (Table1) range, expanding the table automatically.
$ws.Rows.Item(2).Insert()

# Fill the new row 2 with the "vitreous" layer entry.
$ws.Cells.Item(2, 1).Value = 0
$ws.Cells.Item(2, 2).Value = "VIT"
$ws.Cells.Item(2, 3).Value = "vitreous"
$ws.Cells.Item(2, 4).Value = "NA"
$ws.Cells.Item(2, 5).Value = "NA"
$ws.Cells.Item(2, 6).Value = "NA"
$ws.Cells.Item(2, 7).Value = "NA"

# Renumber surface_id values (column A) for the remaining rows so they run
# sequentially 1..11 instead of 1..9,11,12.
$ws.Cells.Item(3, 1).Value = 1
$ws.Cells.Item(4, 1).Value = 2
$ws.Cells.Item(5, 1).Value = 3
$ws.Cells.Item(6, 1).Value = 4
$ws.Cells.Item(7, 1).Value = 5
$ws.Cells.Item(8, 1).Value = 6
$ws.Cells.Item(9, 1).Value = 7
$ws.Cells.Item(10, 1).Value = 8
$ws.Cells.Item(11, 1).Value = 9
$ws.Cells.Item(12, 1).Value = 10
$ws.Cells.Item(13, 1).Value = 11

$ws.Range("A13").Select()

# Resize Table1's range / autofilter to A1:G13 since the new row was inserted
# inside it (covers the case where the insert didn't auto-expand the table).
$tbl = $ws.ListObjects.Item("Table1")
$tbl.Resize($ws.Range("A1:G13"))
